# Update the "想去人数" (want-to-go count) figures for the first two events.
# Both the "展览" sheet and the "全部类型" sheet carry the same data table
# and need to be bumped by +1 each.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 159
    $ws.Range("F3").Value = 110
}
